$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(424).Insert()

$ws.Range("A424").Value = 11
$ws.Range("B424").Value = "Vega Monumental Concepción"
$ws.Range("C424").Value = "Bíobío"
$ws.Range("D424").Value = 45258
$ws.Range("E424").Value = 8
$ws.Range("F424").Value = "Fruta"
$ws.Range("G424").Value = 100103
$ws.Range("H424").Value = "Frutos de hueso (carozo)"
$ws.Range("I424").Value = 100103006
$ws.Range("J424").Value = "Nectarín"
$ws.Range("K424").Value = "Early Glo"
$ws.Range("L424").Value = "Primera"
$ws.Range("M424").Value = 100
$ws.Range("N424").Value = 13000
$ws.Range("O424").Value = 14000
$ws.Range("P424").Value = 13500
$ws.Range("Q424").Value = "$/bandeja 15 kilos granel"
$ws.Range("R424").Value = "Región de O'Higgins"
$ws.Range("S424").Value = 900
$ws.Range("T424").Value = 15
